$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the table; it belongs right before
# the current row 645, so every existing row from 645 down shifts by one.
$ws.Rows(645).Insert()

$ws.Range("A645").Value = 4
$ws.Range("B645").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C645").Value = "Los Lagos"
$ws.Range("D645").Value = 44931
$ws.Range("E645").Value = 10
$ws.Range("F645").Value = "Fruta"
$ws.Range("G645").Value = 100102
$ws.Range("H645").Value = "Cítricos"
$ws.Range("I645").Value = 100102003
$ws.Range("J645").Value = "Limón"
$ws.Range("K645").Value = "Sin especificar"
$ws.Range("L645").Value = "1a plateado"
$ws.Range("M645").Value = 1000
$ws.Range("N645").Value = 21000
$ws.Range("O645").Value = 21000
$ws.Range("P645").Value = 21000
$ws.Range("Q645").Value = "$/malla 18 kilos"
$ws.Range("R645").Value = "Provincia de Melipilla"
$ws.Range("S645").Value = 1167
$ws.Range("T645").Value = 18
